$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'45.755.46"
$ws.Range("E2").Value = "  -2.23%  "

# Row 3
$ws.Range("D3").Value = "'2.450.27"
$ws.Range("E3").Value = "  +8.24%  "

# Row 4
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.23%  "

# Row 5
$ws.Range("D5").Value = "'294.23"
$ws.Range("E5").Value = "  -1.35%  "

# Row 6
$ws.Range("D6").Value = "'95.04"
$ws.Range("E6").Value = "  -2.95%  "

# Row 7
$ws.Range("D7").Value = "'0.570"
$ws.Range("E7").Value = "  -0.82%  "

# Row 8
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("D9").Value = "'0.513"
$ws.Range("E9").Value = "  +2.21%  "

# Row 10
$ws.Range("D10").Value = "'34.92"
$ws.Range("E10").Value = "  +0.68%  "

# Row 11
$ws.Range("D11").Value = "'0.0781"
$ws.Range("E11").Value = "  -1.31%  "

# Row 12
$ws.Range("D12").Value = "'7.24"
$ws.Range("E12").Value = "  +3.71%  "

# Row 13
$ws.Range("E13").Value = "  +1.81%  "

# Row 14
$ws.Range("D14").Value = "'2.828.32"
$ws.Range("E14").Value = "  +8.40%  "

# Row 15
$ws.Range("D15").Value = "'2.445.04"
$ws.Range("E15").Value = "  +7.85%  "

# Row 16
$ws.Range("D16").Value = "'0.849"
$ws.Range("E16").Value = "  +7.37%  "

# Row 17
$ws.Range("D17").Value = "'14.16"
$ws.Range("E17").Value = "  +4.32%  "

# Row 18
$ws.Range("D18").Value = "'45.737.78"
$ws.Range("E18").Value = "  -2.23%  "

# Row 19
$ws.Range("D19").Value = "'12.62"
$ws.Range("E19").Value = "  +2.01%  "

# Row 20
$ws.Range("D20").Value = "'0.0₃0940"
$ws.Range("E20").Value = "  -2.32%  "

# Row 21
$ws.Range("D21").Value = "'6.26"
$ws.Range("E21").Value = "  +8.55%  "

# Row 22
$ws.Range("D22").Value = "'67.48"
$ws.Range("E22").Value = "  +2.70%  "

# Row 23
$ws.Range("D23").Value = "'244.35"
$ws.Range("E23").Value = "  -0.12%  "

# Row 24
$ws.Range("D24").Value = "'2.78"
$ws.Range("E24").Value = "  +0.47%  "

# Row 25
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  -0.05%  "

# Row 26
$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D26").Value = "'1.94"
$ws.Range("E26").Value = "  +5.49%  "

# Row 27
$ws.Range("D27").Value = "'38.99"
$ws.Range("E27").Value = "  -4.96%  "

# Row 28
$ws.Range("D28").Value = "'2.22"
$ws.Range("E28").Value = "  +0.26%  "

# Row 29
$ws.Range("D29").Value = "'9.76"
$ws.Range("E29").Value = "  +2.80%  "

# Row 30
$ws.Range("D30").Value = "'21.54"
$ws.Range("E30").Value = "  +7.57%  "

# Row 31
$ws.Range("D31").Value = "'3.77"
$ws.Range("E31").Value = "  +13.93%  "

# Row 32
$ws.Range("D32").Value = "'2.76"
$ws.Range("E32").Value = "  -2.28%  "

# Row 33
$ws.Range("D33").Value = "'5.52"
$ws.Range("E33").Value = "  +4.71%  "

# Row 34
$ws.Range("D34").Value = "'146.64"
$ws.Range("E34").Value = "  +1.16%  "

# Row 35
$ws.Range("D35").Value = "'2.04"
$ws.Range("E35").Value = "  +24.05%  "

# Row 36
$ws.Range("D36").Value = "'0.0766"
$ws.Range("E36").Value = "  +0.60%  "

# Row 37
$ws.Range("D37").Value = "'0.115"
$ws.Range("E37").Value = "  +3.50%  "

# Row 38
$ws.Range("D38").Value = "'0.115"
$ws.Range("E38").Value = "  +0.05%  "

# Row 39
$ws.Range("D39").Value = "'14.99"
$ws.Range("E39").Value = "  -2.13%  "

# Row 40
$ws.Range("D40").Value = "'3.92"
$ws.Range("E40").Value = "  +3.21%  "

# Row 41
$ws.Range("D41").Value = "'0.0299"
$ws.Range("E41").Value = "  +2.49%  "

# Row 42
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "'2.004.99"
$ws.Range("E42").Value = "  +12.66%  "

# Row 43
$ws.Range("B43").Value = "NEARProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D43").Value = "'3.22"
$ws.Range("E43").Value = "  +5.58%  "

# Row 44
$ws.Range("D44").Value = "'0.998"
$ws.Range("E44").Value = "  -0.17%  "

# Row 45
$ws.Range("D45").Value = "'91.51"
$ws.Range("E45").Value = "  -1.59%  "

# Row 46
$ws.Range("D46").Value = "'16.24"
$ws.Range("E46").Value = "  +32.57%  "

# Row 47
$ws.Range("D47").Value = "'1.77"
$ws.Range("E47").Value = "  -4.98%  "

# Row 48
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'102.58"
$ws.Range("E48").Value = "  +9.61%  "

# Row 49
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "'8.54"
$ws.Range("E49").Value = "  +9.53%  "

# Row 50
$ws.Range("D50").Value = "'2.694.27"
$ws.Range("E50").Value = "  +8.38%  "

# Row 51
$ws.Range("D51").Value = "'0.185"
$ws.Range("E51").Value = "  +1.45%  "
